$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so values like
# "240.64" or "0.624" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '41.912.40'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '2.211.34'
$ws.Range('E3').Value = '  -1.52%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '240.64'
$ws.Range('E5').Value = '  -2.34%  '
$ws.Range('D6').Value = '0.624'
$ws.Range('E6').Value = '  -1.08%  '
$ws.Range('D7').Value = '72.11'
$ws.Range('E7').Value = '  -5.17%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '0.601'
$ws.Range('E9').Value = '  -3.56%  '
$ws.Range('D10').Value = '41.74'
$ws.Range('E10').Value = '  -5.00%  '
$ws.Range('D11').Value = '0.0947'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '6.95'
$ws.Range('E12').Value = '  -4.76%  '
$ws.Range('E13').Value = '  +0.11%  '
$ws.Range('D14').Value = '2.545.57'
$ws.Range('D15').Value = '14.17'
$ws.Range('E15').Value = '  -3.02%  '
$ws.Range('D16').Value = '0.829'
$ws.Range('E16').Value = '  -3.25%  '
$ws.Range('D17').Value = '2.203.63'
$ws.Range('E17').Value = '  -1.79%  '
$ws.Range('D18').Value = '41.809.40'
$ws.Range('E18').Value = '  -0.88%  '
$ws.Range('E19').Value = '  +2.58%  '
$ws.Range('D20').Value = '72.39'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').Value = '6.13'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').Value = '10.89'
$ws.Range('E22').Value = '  +18.65%  '
$ws.Range('D23').Value = '228.97'
$ws.Range('E23').Value = '  -1.15%  '
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -8.67%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').Value = '11.42'
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('E28').Value = '  -2.52%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').Value = '167.30'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').Value = '20.40'
$ws.Range('E31').Value = '  -1.48%  '
$ws.Range('D32').Value = '5.57'
$ws.Range('E32').Value = '  +5.58%  '
$ws.Range('D33').Value = '0.0791'
$ws.Range('E33').Value = '  -4.64%  '
$ws.Range('D34').Value = '29.97'
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('E36').Value = '  -11.95%  '
$ws.Range('D37').Value = '4.20'
$ws.Range('E37').Value = '  -7.36%  '
$ws.Range('D38').Value = '0.0298'
$ws.Range('E38').Value = '  -6.17%  '
$ws.Range('D39').Value = '13.57'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('E40').Value = '  -3.77%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '5.60'
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').Value = '63.89'
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('D44').Value = '8.65'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = '103.24'
$ws.Range('E45').Value = '  -4.59%  '
$ws.Range('E46').Value = '  -1.89%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '1.10'
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').Value = '1.16'
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '2.31'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').Value = '2.70'
$ws.Range('E50').Value = '  +0.08%  '
$ws.Range('D51').Value = '2.420.22'
$ws.Range('E51').Value = '  -1.49%  '
